$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 11, shifting existing rows (11+) down by one.
$ws.Rows("11:11").Insert()

# Fill in the new paper entry in row 11.
$ws.Range("A11").Value = "Offshore wind installation: Analysing the evidence behind improvements ininstallation time"
$ws.Range("B11").Value = "Lacal-Arantegui et al"
$ws.Range("C11").Value = 2018
$ws.Range("E11").Value = "Short (9 pages, lots of graphs and pics) paper analysing reducement in installation time based on 'better' turbines?"
$ws.Range("F11").Value = "N"
$ws.Range("G11").Value = "M"
$ws.Range("H11").Value = "https://reader.elsevier.com/reader/sd/pii/S1364032118302612?token=11096C5E4A1FEBD54590D279F2295F15FC8ACD7E0CD4715D62E826B87DF809B7427D43E68891E5CEA053E032B4A1E2C9"

$ws.Hyperlinks.Add($ws.Range("H11"), "https://reader.elsevier.com/reader/sd/pii/S1364032118302612?token=11096C5E4A1FEBD54590D279F2295F15FC8ACD7E0CD4715D62E826B87DF809B7427D43E68891E5CEA053E032B4A1E2C9")

$ws.Rows("11:11").RowHeight = 30
